# Dockerized selenium grid setup
$wb = $excel.ActiveWorkbook

$runner = $wb.Worksheets.Item("Runner")
$cred = $wb.Worksheets.Item("CredentialData")

# --- CredentialData sheet: fix casing "Yes" -> "yes" for D3 and D4 ---
# Copy the existing lowercase "yes" text from Runner!C3 first (before it gets
# overwritten below) so the text + cell formatting stay consistent with the
# shared-string table, then paste values only (keep each cell's own style).
$runner.Range("C3").Copy()
$cred.Range("D3").PasteSpecial(-4163)  # xlPasteValues
$cred.Range("D4").PasteSpecial(-4163)  # xlPasteValues

# --- Runner sheet: fix casing of "yes" -> "Yes" in C3 ---
# Copy a cell already containing "Yes" so the cell style (incl. quotePrefix) is preserved.
$runner.Range("C2").Copy()
$runner.Range("C3").PasteSpecial(-4163)  # xlPasteValues

# Apply the style used by body rows on Runner sheet to body rows on CredentialData sheet
# (matches removal of the now-duplicate/unused cellXfs entry)
$runner.Range("A2").Copy()
$cred.Range("A2:F4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the saved selection on CredentialData (A2:F4) ...
$cred.Activate()
$cred.Range("A2:F4").Select()

# ... then restore Runner as the active sheet/selection (C2), matching tabSelected
$runner.Activate()
$runner.Range("C2").Select()
